# Update odds values in row 3 and row 4 of the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 changes
$ws.Range("G3").Value = 1.9
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 2.63
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("X3").Value = 7.5
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 15
$ws.Range("AE3").Value = 21
$ws.Range("AF3").Value = 81
$ws.Range("AG3").Value = 9.5
$ws.Range("AH3").Value = 21
$ws.Range("AJ3").Value = 51
$ws.Range("AU3").Value = 81

# Row 4 changes
$ws.Range("G4").Value = 2.63
$ws.Range("H4").Value = 2.75
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5
$ws.Range("Z4").Value = 29
$ws.Range("AN4").Value = 19
$ws.Range("AT4").Value = 11
$ws.Range("AZ4").Value = 151
